# Apply the targeted updates to the "Work Report" sheet:
#  - D5: update the "Report Generated On" timestamp
#  - C8: Total Billed Amount -> 0
#  - H16: line-item Pricing -> 0
#  - H17: TOTAL Pricing -> 0

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Report")

$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"

$ws.Range("C8").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
